$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.486.05'
$ws.Range("E2").Value = '  +2.06%  '

$ws.Range("D3").Value = '2.674.72'
$ws.Range("E3").Value = '  +2.22%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.09%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.87'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.74%  '

$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.525'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.07%  '

$ws.Range("D9").Value = '2.675.11'
$ws.Range("E9").Value = '  +2.26%  '

$ws.Range("E10").Value = '  +4.68%  '

$ws.Range("E11").Value = '  +2.28%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.358'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.16%  '

$ws.Range("E13").Value = '  +0.28%  '

$ws.Range("D14").Value = '3.164.46'
$ws.Range("E14").Value = '  +2.15%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000187'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.54%  '

$ws.Range("D16").Value = '72.381.02'
$ws.Range("E16").Value = '  +2.00%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.60'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.06%  '

$ws.Range("D18").Value = '2.673.45'
$ws.Range("E18").Value = '  +3.49%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.96'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.44%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.03'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.24%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '380.46'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.08%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.21'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.48%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.07'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +11.91%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.57'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.10%  '

$ws.Range("E25").Value = '  -0.09%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.38'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.71%  '

$ws.Range("E27").Value = '  +3.74%  '

$ws.Range("D28").Value = '2.807.12'
$ws.Range("E28").Value = '  +3.50%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.52%  '

$ws.Range("D30").Value = '0.0₃0951'
$ws.Range("E30").Value = '  -0.23%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.15'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.69%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '522.36'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.12%  '

$ws.Range("E33").Value = '  +0.00%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.83'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.39%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.01%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '163.46'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.33%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.61'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.30%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.14'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.99%  '

$ws.Range("E39").Value = '  +1.57%  '

$ws.Range("E40").Value = '  -6.86%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.85'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.16%  '

$ws.Range("E42").Value = '  -0.04%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.07'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.68%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.60'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.28%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.336'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.39%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.32'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.88%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '153.22'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.03%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.76'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.22%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.551'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.73%  '

$ws.Range("E50").Value = '  +2.74%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0769'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.03%  '
